# Update the "Exchange Conversion Date/Time" column (K) values from
# 05/08/2025 to 06/08/2025 for all data rows (rows 2-10) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 11)  # Column K = 11
    if ($cell.Value2 -eq "05/08/2025") {
        # Prefix with an apostrophe so Excel stores this as literal text
        # instead of auto-converting the date-looking string into a real
        # date value (keeps it matching the original inline/shared string).
        $cell.Value = "'06/08/2025"
    }
}
